# Updated cryptos list with refreshed Price / Volume(1h) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold formatted text values (e.g. thousands
# separated with dots, or percentages padded with spaces). Force text format so
# Excel does not reinterpret them as numbers and strip formatting (trailing zeros,
# padding, etc.).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.036.78'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '1.833.28'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('D4').Value = '0.9985'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '241.64'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').Value = '0.6279'
$ws.Range('E6').Value = '  -5.15%  '
$ws.Range('D7').Value = '1.0000'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.07600'
$ws.Range('E8').Value = '  +2.12%  '
$ws.Range('D9').Value = '0.2915'
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('D10').Value = '22.72'
$ws.Range('E10').Value = '  -2.64%  '
$ws.Range('D11').Value = '0.07737'
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').Value = '1.837.44'
$ws.Range('D13').Value = '4.955'
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('D14').Value = '0.6641'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('D15').Value = '82.74'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').Value = '0.000009803'
$ws.Range('E16').Value = '  +14.19%  '
$ws.Range('D17').Value = '5.988'
$ws.Range('E17').Value = '  -3.06%  '
$ws.Range('D18').Value = '29.027.58'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = '226.57'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('E20').Value = '  -1.74%  '
$ws.Range('D21').Value = '0.9994'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').Value = '7.216'
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '158.22'
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').Value = '0.1372'
$ws.Range('E25').Value = '  -2.31%  '
$ws.Range('D26').Value = '8.415'
$ws.Range('E26').Value = '  -2.50%  '
$ws.Range('E27').Value = '  -1.09%  '
$ws.Range('D28').Value = '1.489'
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('D29').Value = '4.062'
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('D30').Value = '4.023'
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('D31').Value = '1.195'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').Value = '0.05187'
$ws.Range('E32').Value = '  -2.64%  '
$ws.Range('D33').Value = '1.847'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').Value = '0.7393'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('D36').Value = '2.696'
$ws.Range('E36').Value = '  +1.58%  '
$ws.Range('D37').Value = '1.264.74'
$ws.Range('E37').Value = '  -3.96%  '
$ws.Range('D38').Value = '2.759'
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').Value = '0.01787'
$ws.Range('E39').Value = '  -0.79%  '
$ws.Range('D40').Value = '6.254'
$ws.Range('E40').Value = '  -2.54%  '
$ws.Range('D41').Value = '0.8948'
$ws.Range('E41').Value = '  -2.50%  '
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').Value = '101.50'
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('D44').Value = '1.977.50'
$ws.Range('E44').Value = '  -1.77%  '
$ws.Range('D45').Value = '0.00000000124'
$ws.Range('E45').Value = '  +2.46%  '
$ws.Range('D46').Value = '64.56'
$ws.Range('E46').Value = '  -2.04%  '
$ws.Range('D47').Value = '0.5109'
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('D48').Value = '0.3987'
$ws.Range('E48').Value = '  -1.01%  '
$ws.Range('D49').Value = '8.860'
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('D50').Value = '0.05753'
$ws.Range('E50').Value = '  -1.98%  '
$ws.Range('D51').Value = '6.678'
$ws.Range('E51').Value = '  -2.11%  '
